$d = $word.ActiveDocument

# Locate the "UMenu" paragraph under C++ objects > Multiplayer Plugin > UGameInstanceSubsystem
# (the one whose paragraph properties include w:spacing before/after 0, as in the diff context).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "UMenu" -and $p.Range.ListFormat.ListLevelNumber -eq 4) {
        $target = $p
        break
    }
}

$r = $target.Range
$r.InsertParagraphAfter()
$p1 = $target.Next()
$p1.Range.Text = "Blaster"
$p1.Range.ListFormat.ListLevelNumber = 2

$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.Text = "ACharacter"
$p2.Range.ListFormat.ListLevelNumber = 3

$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
$p3.Range.Text = "ABlasterCharacter"
$p3.Range.ListFormat.ListLevelNumber = 4
